# Update inventory exposure simulation with enhanced analysis and presentation plots
# Apply updated Price Impact / Incremental IL / IL-per-Price-Impact values
# for rows 3-11 (columns F, G, H) on the LP Analysis sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = -32.87671232876712
$ws.Range("G3").Value = -30.02005979290482
$ws.Range("H3").Value = 91.31101520341883

$ws.Range("F4").Value = -24.74226804123711
$ws.Range("G4").Value = -19.38281475247133
$ws.Range("H4").Value = 78.33887629123831

$ws.Range("F5").Value = -19.83471074380164
$ws.Range("G5").Value = -14.00493907779026
$ws.Range("H5").Value = 70.60823451719259

$ws.Range("F6").Value = -16.55172413793105
$ws.Range("G6").Value = -10.79937416041439
$ws.Range("H6").Value = 65.24621888583688

$ws.Range("F7").Value = -14.20118343195265
$ws.Range("G7").Value = -8.687173283173244
$ws.Range("H7").Value = 61.17217853567833

$ws.Range("F8").Value = -12.43523316062176
$ws.Range("G8").Value = -7.197972117526796
$ws.Range("H8").Value = 57.88369244511132

$ws.Range("F9").Value = -11.05990783410138
$ws.Range("G9").Value = -6.095481676379199
$ws.Range("H9").Value = 55.11331349059527

$ws.Range("F10").Value = -9.958506224066399
$ws.Range("G10").Value = -5.248584357166031
$ws.Range("H10").Value = 52.70453458654218

$ws.Range("F11").Value = 0
$ws.Range("G11").Value = 0
$ws.Range("H11").Value = 0
